$wb = $excel.ActiveWorkbook

# Add the new "Campaign" worksheet as the last tab (after "Activity"),
# mirroring how the sheet was inserted at the end of the sheets collection.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$campaign = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$campaign.Name = "Campaign"

# Populate the two rows: bold header "RecordType" and data value "Parent Campaign".
$campaign.Range("A1").Value = "RecordType"
$campaign.Range("A1").Font.Bold = $true
$campaign.Range("A2").Value = "Parent Campaign"

# Match the column A "best fit" width used for the header/value text.
$campaign.Columns.Item(1).ColumnWidth = 14

# Make the newly added Campaign sheet the active/selected tab.
$campaign.Activate()
